# Update "想去人数" (interest count) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 52
    6  = 476
    7  = 1337
    8  = 458
    9  = 91
    10 = 162
    11 = 114
    12 = 163
    14 = 141
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
